$wb = $excel.ActiveWorkbook

# --- Sheet 1 "table attribute": drop the "attribute value" column, keep a
# single column with the table name followed by its value. ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B1:B2").Clear()
$ws1.Range("A1").Value = "table name "
$ws1.Range("A2").Value = "log packet struction"
[void]$ws1.Range("A5").Select()

# --- Sheet 2 "table content Columns": rename the value-type column and add
# a new "column description" column, plus a page setup. ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B1").Value = "column type"
$ws2.Range("C1").Value = "column description"
$ws2.Columns.Item(3).ColumnWidth = 21.65
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# --- Sheet 3 "table content": just loses the "active" selection, no data
# changes. ---
$ws3 = $wb.Worksheets.Item(3)
[void]$ws3.Range("E15").Select()

# Sheet 2 becomes the active/selected tab (activeTab goes from 2 to 1).
[void]$ws2.Range("D6").Select()
$ws2.Activate()
